$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D and E columns (rows 2-51) before writing values so
# numeric-looking strings (e.g. "1.000", "0.9998") are preserved as text,
# matching the original inline-string cell type.
$priceRange = $ws.Range("D2:D51")
$volRange = $ws.Range("E2:E51")
$priceRange.NumberFormat = "@"
$volRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.397.56"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.848.42"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("D5").Value = "233.32"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.4674"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("D8").Value = "0.2729"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "0.06297"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("D10").Value = "1.829.64"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").Value = "0.07457"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "16.35"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "4.928"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "83.87"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "0.6209"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "30.333.62"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D18").Value = "229.02"
$ws.Range("E18").Value = "  +1.52%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "12.38"
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "4.909"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("D23").Value = "5.875"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("D24").Value = "166.80"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "17.82"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "1.871"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "1.376"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "4.084"
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("D31").Value = "3.812"
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").Value = "0.04882"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "1.138"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").Value = "0.7012"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("D35").Value = "2.689"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "0.01920"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "2.655"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").Value = "0.8655"
$ws.Range("E38").Value = "  -3.94%  "
$ws.Range("D39").Value = "105.85"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "1.939"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "5.514"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "0.4031"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").Value = "7.072"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").Value = "61.04"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").Value = "0.1204"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "8.562"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").Value = "33.36"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("D49").Value = "0.05543"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").Value = "0.3646"
$ws.Range("E51").Value = "  -1.57%  "

# Clear the temporary text-number-format styling so the cells go back to
# the default (unstyled) cell format, matching the original workbook.
$priceRange.ClearFormats()
$volRange.ClearFormats()

